$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly observation was added to the dataset: insert a new row
# at position 366, pushing every following row (old 366..392) down by
# one (they become 367..393). This also grows the used range from
# A1:R392 to A1:R393.
$ws.Rows.Item(366).Insert()

# Populate the newly inserted row 366 with the new record's data.
$ws.Range("A366").Value = 3
$ws.Range("B366").Value = "Femacal de La Calera"
$ws.Range("C366").Value = "Coquimbo"
$ws.Range("D366").Value = 44746
$ws.Range("E366").Value = 5
$ws.Range("F366").Value = 100112031
$ws.Range("G366").Value = "Poroto verde"
$ws.Range("H366").Value = "Magnum"
$ws.Range("I366").Value = "Primera"
$ws.Range("J366").Value = 80
$ws.Range("K366").Value = 33000
$ws.Range("L366").Value = 34000
$ws.Range("M366").Value = 33500
$ws.Range("N366").Value = "$/malla 25 kilos"
$ws.Range("O366").Value = "Región de Arica y Parinacota"
$ws.Range("P366").Value = 1340
$ws.Range("Q366").Value = 25
$ws.Range("R366").Value = "Hortaliza"
